# Atualização automática de IRAI.xlsx
#
# This script reproduces the author's edit:
#   1. Rename "Paineis DARQ"            -> "PAINEIS DARQ"
#   2. Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
#   3. Delete the "Desarquivamentos Pendentes" worksheet entirely
#
# (The "DGC" sheet itself keeps its original content; it only appears to
# change in the raw OOXML because deleting "Desarquivamentos Pendentes"
# shifts the shared-string table / sheet numbering that DGC's cells and
# the workbook's sheetId counters point into - that bookkeeping is handled
# automatically by Excel when the sheet is removed.)

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()
$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"
